$d = $word.ActiveDocument

# The paragraph "${obs_carpeta} ${negrita_carpeta} " (unique in the
# document) has a stray trailing space after the closing "}" of
# ${negrita_carpeta}. That trailing space is its own bold run; removing
# it is exactly the edit described by the diff (the run is deleted
# outright, nothing else in the paragraph changes).
$needle = "negrita_carpeta} "
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains($needle)) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    # $r.End sits just past the paragraph mark, so the character right
    # before the mark ($r.End - 2 .. $r.End - 1) is the trailing space.
    $spaceRange = $d.Range($r.End - 2, $r.End - 1)
    if ($spaceRange.Text -eq " ") {
        $spaceRange.Delete()
    }
}
